$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated "K" column (G) values replacing the old "Strike#" counts.
# Maps row -> new K value, per the recalculated save_data for hall_dl.xlsx
$newK = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 1
    6  = 4
    7  = 1
    8  = 6
    9  = 5
    10 = 9
    11 = 3
    12 = 4
    13 = 5
    14 = 1
    15 = 4
    16 = 3
    17 = 3
    18 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
